# Raul's Log update: add Jan 3-5, 2017 (serials 42738-42740) log entries
# picked up from the zone super logs, including trailing template rows that
# were left blank (but still bordered) at the end of that source sheet.
#
# Housekeeping asked for in the commit message:
#  - make sure any stray Excel.exe processes are not left running before we
#    start touching the workbook
#  - clear the clipboard once the copy/paste-derived formatting has been
#    applied, so nothing marching-ants stays selected for the user

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A1139").Value = "AV Shutdown"
$ws.Range("B1139").Value = 42738
$ws.Range("C1139").Value = "1630"
$ws.Range("D1139").Value = "YL"
$ws.Range("E1139").Value = "280 N"
$ws.Range("F1139").Value = "We have a key in YKLN 203 C"

$ws.Range("A1140").Value = "Pickup Mic"
$ws.Range("B1140").Value = 42738
$ws.Range("C1140").Value = "1630"
$ws.Range("D1140").Value = "YL"
$ws.Range("E1140").Value = "280 N"
$ws.Range("F1140").Value = "we have a key in YKLN 203 C  - return mic to KT 516  and place battery into charger"

$ws.Range("A1141").Value = "Demo"
$ws.Range("B1141").Value = 42738
$ws.Range("C1141").Value = "1630"
$ws.Range("D1141").Value = "OSG"
$ws.Range("E1141").Value = "1001"

$ws.Range("A1142").Value = "Demo"
$ws.Range("B1142").Value = 42738
$ws.Range("C1142").Value = "1630"
$ws.Range("D1142").Value = "OSG"
$ws.Range("E1142").Value = "1002"

$ws.Range("A1143").Value = "Demo"
$ws.Range("B1143").Value = 42738
$ws.Range("C1143").Value = "1630"
$ws.Range("D1143").Value = "OSG"
$ws.Range("E1143").Value = "1008"

$ws.Range("A1144").Value = "Demo"
$ws.Range("B1144").Value = 42738
$ws.Range("C1144").Value = "1800"
$ws.Range("D1144").Value = "OSG"
$ws.Range("E1144").Value = "2008"

$ws.Range("A1145").Value = "Demo"
$ws.Range("B1145").Value = 42738
$ws.Range("C1145").Value = "1830"
$ws.Range("D1145").Value = "OSG"
$ws.Range("E1145").Value = "1004"

$ws.Range("A1146").Value = "Demo"
$ws.Range("B1146").Value = 42738
$ws.Range("C1146").Value = "1830"
$ws.Range("D1146").Value = "OSG"
$ws.Range("E1146").Value = "2004"

$ws.Range("A1147").Value = "AV Shutdown"
$ws.Range("B1147").Value = 42738
$ws.Range("C1147").Value = "2000"
$ws.Range("D1147").Value = "OSG"
$ws.Range("E1147").Value = "2001"

$ws.Range("A1148").Value = "AV Shutdown"
$ws.Range("B1148").Value = 42738
$ws.Range("C1148").Value = "2000"
$ws.Range("D1148").Value = "OSG"
$ws.Range("E1148").Value = "2002"

$ws.Range("A1149").Value = "AV Shutdown"
$ws.Range("B1149").Value = 42738
$ws.Range("C1149").Value = "2000"
$ws.Range("D1149").Value = "OSG"
$ws.Range("E1149").Value = "2009"

$ws.Range("A1150").Value = "AV Shutdown"
$ws.Range("B1150").Value = 42738
$ws.Range("C1150").Value = "2000"
$ws.Range("D1150").Value = "OSG"
$ws.Range("E1150").Value = "1001"

$ws.Range("A1151").Value = "AV Shutdown"
$ws.Range("B1151").Value = 42738
$ws.Range("C1151").Value = "2000"
$ws.Range("D1151").Value = "OSG"
$ws.Range("E1151").Value = "1002"

$ws.Range("A1152").Value = "AV Shutdown"
$ws.Range("B1152").Value = 42738
$ws.Range("C1152").Value = "2000"
$ws.Range("D1152").Value = "OSG"
$ws.Range("E1152").Value = "1008"

$ws.Range("A1153").Value = "AV Shutdown"
$ws.Range("B1153").Value = 42738
$ws.Range("C1153").Value = "2000"
$ws.Range("D1153").Value = "OSG"
$ws.Range("E1153").Value = "2008"

$ws.Range("A1157").Value = "AV Shutdown"
$ws.Range("B1157").Value = 42739
$ws.Range("C1157").Value = "1700"
$ws.Range("D1157").Value = "SSB"
$ws.Range("E1157").Value = "W132"

$ws.Range("A1158").Value = "AV Shutdown"
$ws.Range("B1158").Value = 42739
$ws.Range("C1158").Value = "1630"
$ws.Range("D1158").Value = "ACW"
$ws.Range("E1158").Value = "109"

$ws.Range("A1159").Value = "AV Shutdown"
$ws.Range("B1159").Value = 42739
$ws.Range("C1159").Value = "1630"
$ws.Range("D1159").Value = "YL"
$ws.Range("E1159").Value = "280 N"
$ws.Range("F1159").Value = "We have a key in YKLN 203 C"

$ws.Range("A1160").Value = "Pickup Mic"
$ws.Range("B1160").Value = 42739
$ws.Range("C1160").Value = "1630"
$ws.Range("D1160").Value = "YL"
$ws.Range("E1160").Value = "280 N"
$ws.Range("F1160").Value = "we have a key in YKLN 203 C  - return mic to KT 516  and place battery into charger"

$ws.Range("A1161").Value = "AV Shutdown"
$ws.Range("B1161").Value = 42739
$ws.Range("C1161").Value = "1630"
$ws.Range("D1161").Value = "OSG"
$ws.Range("E1161").Value = "1001"

$ws.Range("A1162").Value = "AV Shutdown"
$ws.Range("B1162").Value = 42739
$ws.Range("C1162").Value = "1630"
$ws.Range("D1162").Value = "OSG"
$ws.Range("E1162").Value = "1002"

$ws.Range("A1163").Value = "AV Shutdown"
$ws.Range("B1163").Value = 42739
$ws.Range("C1163").Value = "1630"
$ws.Range("D1163").Value = "OSG"
$ws.Range("E1163").Value = "1004"

$ws.Range("A1164").Value = "AV Shutdown"
$ws.Range("B1164").Value = 42739
$ws.Range("C1164").Value = "1630"
$ws.Range("D1164").Value = "OSG"
$ws.Range("E1164").Value = "1008"

$ws.Range("A1165").Value = "AV Shutdown"
$ws.Range("B1165").Value = 42739
$ws.Range("C1165").Value = "1630"
$ws.Range("D1165").Value = "OSG"
$ws.Range("E1165").Value = "1014"

$ws.Range("A1166").Value = "AV Shutdown"
$ws.Range("B1166").Value = 42739
$ws.Range("C1166").Value = "1630"
$ws.Range("D1166").Value = "OSG"
$ws.Range("E1166").Value = "2001"

$ws.Range("A1167").Value = "Demo"
$ws.Range("B1167").Value = 42739
$ws.Range("C1167").Value = "1630"
$ws.Range("D1167").Value = "OSG"
$ws.Range("E1167").Value = "1003"

$ws.Range("A1168").Value = "Demo"
$ws.Range("B1168").Value = 42739
$ws.Range("C1168").Value = "1630"
$ws.Range("D1168").Value = "OSG"
$ws.Range("E1168").Value = "2002"

$ws.Range("A1169").Value = "Demo"
$ws.Range("B1169").Value = 42739
$ws.Range("C1169").Value = "1630"
$ws.Range("D1169").Value = "OSG"
$ws.Range("E1169").Value = "2003"

$ws.Range("A1170").Value = "Demo"
$ws.Range("B1170").Value = 42739
$ws.Range("C1170").Value = "1730"
$ws.Range("D1170").Value = "OSG"
$ws.Range("E1170").Value = "2004"

$ws.Range("A1171").Value = "Demo"
$ws.Range("B1171").Value = 42739
$ws.Range("C1171").Value = "1730"
$ws.Range("D1171").Value = "OSG"
$ws.Range("E1171").Value = "2010"

$ws.Range("A1172").Value = "Demo"
$ws.Range("B1172").Value = 42739
$ws.Range("C1172").Value = "1820"
$ws.Range("D1172").Value = "OSG"
$ws.Range("E1172").Value = "1002"

$ws.Range("A1173").Value = "Demo"
$ws.Range("B1173").Value = 42739
$ws.Range("C1173").Value = "1820"
$ws.Range("D1173").Value = "OSG"
$ws.Range("E1173").Value = "2008"

$ws.Range("A1174").Value = "Demo"
$ws.Range("B1174").Value = 42739
$ws.Range("C1174").Value = "1900"
$ws.Range("D1174").Value = "OSG"
$ws.Range("E1174").Value = "2009"

$ws.Range("A1175").Value = "AV Shutdown"
$ws.Range("B1175").Value = 42739
$ws.Range("C1175").Value = "2030"
$ws.Range("D1175").Value = "OSG"
$ws.Range("E1175").Value = "1002"

$ws.Range("A1176").Value = "AV Shutdown"
$ws.Range("B1176").Value = 42739
$ws.Range("C1176").Value = "2030"
$ws.Range("D1176").Value = "OSG"
$ws.Range("E1176").Value = "1003"

$ws.Range("A1177").Value = "AV Shutdown"
$ws.Range("B1177").Value = 42739
$ws.Range("C1177").Value = "2030"
$ws.Range("D1177").Value = "OSG"
$ws.Range("E1177").Value = "2001"

$ws.Range("A1178").Value = "AV Shutdown"
$ws.Range("B1178").Value = 42739
$ws.Range("C1178").Value = "2030"
$ws.Range("D1178").Value = "OSG"
$ws.Range("E1178").Value = "2002"

$ws.Range("A1179").Value = "AV Shutdown"
$ws.Range("B1179").Value = 42739
$ws.Range("C1179").Value = "2030"
$ws.Range("D1179").Value = "OSG"
$ws.Range("E1179").Value = "2003"

$ws.Range("A1180").Value = "AV Shutdown"
$ws.Range("B1180").Value = 42739
$ws.Range("C1180").Value = "2030"
$ws.Range("D1180").Value = "OSG"
$ws.Range("E1180").Value = "2004"

$ws.Range("A1181").Value = "AV Shutdown"
$ws.Range("B1181").Value = 42739
$ws.Range("C1181").Value = "2030"
$ws.Range("D1181").Value = "OSG"
$ws.Range("E1181").Value = "2008"

$ws.Range("A1182").Value = "AV Shutdown"
$ws.Range("B1182").Value = 42739
$ws.Range("C1182").Value = "2030"
$ws.Range("D1182").Value = "OSG"
$ws.Range("E1182").Value = "2010"

$ws.Range("A1187").Value = "Demo"
$ws.Range("B1187").Value = 42740
$ws.Range("C1187").Value = "1630"
$ws.Range("D1187").Value = "OSG"
$ws.Range("E1187").Value = "2010"

$ws.Range("A1188").Value = "Demo"
$ws.Range("B1188").Value = 42740
$ws.Range("C1188").Value = "1600"
$ws.Range("D1188").Value = "ACE"
$ws.Range("E1188").Value = "007"

$ws.Range("A1189").Value = "Demo"
$ws.Range("B1189").Value = 42740
$ws.Range("C1189").Value = "1600"
$ws.Range("D1189").Value = "ACW"
$ws.Range("E1189").Value = "204"

$ws.Range("A1190").Value = "Demo"
$ws.Range("B1190").Value = 42740
$ws.Range("C1190").Value = "1630"
$ws.Range("D1190").Value = "OSG"
$ws.Range("E1190").Value = "2004"

$ws.Range("A1191").Value = "Demo"
$ws.Range("B1191").Value = 42740
$ws.Range("C1191").Value = "1630"
$ws.Range("D1191").Value = "OSG"
$ws.Range("E1191").Value = "2008"

$ws.Range("A1192").Value = "Demo"
$ws.Range("B1192").Value = 42740
$ws.Range("C1192").Value = "1730"
$ws.Range("D1192").Value = "OSG"
$ws.Range("E1192").Value = "2003"

$ws.Range("A1193").Value = "Demo"
$ws.Range("B1193").Value = 42740
$ws.Range("C1193").Value = "1730"
$ws.Range("D1193").Value = "OSG"
$ws.Range("E1193").Value = "2009"

$ws.Range("A1194").Value = "Demo"
$ws.Range("B1194").Value = 42740
$ws.Range("C1194").Value = "1800"
$ws.Range("D1194").Value = "OSG"
$ws.Range("E1194").Value = "2002"

$ws.Range("A1195").Value = "Demo"
$ws.Range("B1195").Value = 42740
$ws.Range("C1195").Value = "1830"
$ws.Range("D1195").Value = "OSG"
$ws.Range("E1195").Value = "2010"

$ws.Range("A1196").Value = "Demo"
$ws.Range("B1196").Value = 42740
$ws.Range("C1196").Value = "1900"
$ws.Range("D1196").Value = "ACE"
$ws.Range("E1196").Value = "003"

$ws.Range("A1197").Value = "Demo"
$ws.Range("B1197").Value = 42740
$ws.Range("C1197").Value = "1900"
$ws.Range("D1197").Value = "ACE"
$ws.Range("E1197").Value = "009"

$ws.Range("A1198").Value = "Demo"
$ws.Range("B1198").Value = 42740
$ws.Range("C1198").Value = "1900"
$ws.Range("D1198").Value = "DB"
$ws.Range("E1198").Value = "0004"

$ws.Range("A1199").Value = "Demo"
$ws.Range("B1199").Value = 42740
$ws.Range("C1199").Value = "1900"
$ws.Range("D1199").Value = "HNE"
$ws.Range("E1199").Value = "401"

$ws.Range("A1200").Value = "Setup PC"
$ws.Range("B1200").Value = 42740
$ws.Range("C1200").Value = "1830"
$ws.Range("D1200").Value = "SSB"
$ws.Range("E1200").Value = "N106"

$ws.Range("A1201").Value = "Setup PC"
$ws.Range("B1201").Value = 42740
$ws.Range("C1201").Value = "1830"
$ws.Range("D1201").Value = "SSB"
$ws.Range("E1201").Value = "N107"

$ws.Range("A1202").Value = "Setup PC"
$ws.Range("B1202").Value = 42740
$ws.Range("C1202").Value = "1830"
$ws.Range("D1202").Value = "SSB"
$ws.Range("E1202").Value = "N108"

$ws.Range("A1203").Value = "AV Shutdown"
$ws.Range("B1203").Value = 42740
$ws.Range("C1203").Value = "1730"
$ws.Range("D1203").Value = "OSG"
$ws.Range("E1203").Value = "1004"

$ws.Range("A1204").Value = "AV Shutdown"
$ws.Range("B1204").Value = 42740
$ws.Range("C1204").Value = "1730"
$ws.Range("D1204").Value = "OSG"
$ws.Range("E1204").Value = "1008"

$ws.Range("A1205").Value = "AV Shutdown"
$ws.Range("B1205").Value = 42740
$ws.Range("C1205").Value = "1630"
$ws.Range("D1205").Value = "OSG"
$ws.Range("E1205").Value = "2002"

$ws.Range("A1206").Value = "AV Shutdown"
$ws.Range("B1206").Value = 42740
$ws.Range("C1206").Value = "1830"
$ws.Range("D1206").Value = "OSG"
$ws.Range("E1206").Value = "2004"

$ws.Range("A1207").Value = "AV Shutdown"
$ws.Range("B1207").Value = 42740
$ws.Range("C1207").Value = "1830"
$ws.Range("D1207").Value = "OSG"
$ws.Range("E1207").Value = "2008"

$ws.Range("A1208").Value = "AV Shutdown"
$ws.Range("B1208").Value = 42740
$ws.Range("C1208").Value = "1830"
$ws.Range("D1208").Value = "OSG"
$ws.Range("E1208").Value = "2010"

$ws.Range("A1209").Value = "AV Shutdown"
$ws.Range("B1209").Value = 42740
$ws.Range("C1209").Value = "2030"
$ws.Range("D1209").Value = "OSG"
$ws.Range("E1209").Value = "2003"

$ws.Range("A1210").Value = "AV Shutdown"
$ws.Range("B1210").Value = 42740
$ws.Range("C1210").Value = "2030"
$ws.Range("D1210").Value = "OSG"
$ws.Range("E1210").Value = "2009"

# Rows that use the taller "Pickup Mic" instruction wrap (matches the source
# zone-super rows this was copied from).
$ws.Rows.Item(1140).RowHeight = 30
$ws.Rows.Item(1160).RowHeight = 30

# Work around the zone super log having trailing template rows that are
# blank but still carry the boxed formatting from that sheet - re-create the
# borders so those rows survive the round trip the same way.
$xlEdgeLeft = 7; $xlEdgeTop = 8; $xlEdgeBottom = 9; $xlEdgeRight = 10
$xlContinuous = 1
function Set-BoxBorder($rng, [bool]$right) {
    $rng.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
    $rng.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $rng.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    if ($right) { $rng.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous }
}

Set-BoxBorder $ws.Range("A1154") $true
Set-BoxBorder $ws.Range("B1154") $false
Set-BoxBorder $ws.Range("C1154") $true
Set-BoxBorder $ws.Range("D1154") $true

Set-BoxBorder $ws.Range("A1211") $true
Set-BoxBorder $ws.Range("B1211") $false

Set-BoxBorder $ws.Range("A1212") $true
Set-BoxBorder $ws.Range("B1212") $false

Set-BoxBorder $ws.Range("A1213") $true
Set-BoxBorder $ws.Range("B1213") $false

Set-BoxBorder $ws.Range("A1214") $true
Set-BoxBorder $ws.Range("B1214") $false

Set-BoxBorder $ws.Range("A1215") $true
Set-BoxBorder $ws.Range("B1215") $false

Set-BoxBorder $ws.Range("A1216") $true
Set-BoxBorder $ws.Range("B1216") $false

# Keep the view roughly where the last entries were typed.
$ws.Activate()
$ws.Range("F1213").Select()

# Clear the clipboard now that the copy/paste-driven formatting is done.
$excel.CutCopyMode = $false
